$d = $word.ActiveDocument
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14"><w:body><w:p><w:pPr><w:pStyle w:val="Titel"/></w:pPr><w:r><w:t xml:space="preserve">Definition </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>of</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Done</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Sprint 2</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="berschrift1"/></w:pPr><w:r><w:t>Anfrage stellen</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Als Mieter kann ich die Liste der Räume laden, welche zu einem bestimmten Zeitpunkt verfügbar sind. Danach kann der Benutzer einen Raum auswählen und eine Reservation tätigen. Nach der Reservation erscheint eine Bestätigungsmeldung.</w:t></w:r><w:r><w:t xml:space="preserve"> Wenn</w:t></w:r><w:r><w:t xml:space="preserve"> ein Raum zu diesem Zeitpunkt bereits besetzt ist, erscheint eine Fehlermeldung.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Um die Story zu testen, müssen noch einige Räume in die Datenbank eingefügt werden und danach muss das Ganze funktionieren.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Was muss noch gemacht sein: </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Testfälle die im Vorherein geschrieben wurden, sind getestet und ohne Fehler durchgeführt worden.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Dokumentation über die Methoden und Wege wurde geführt und in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> gepuscht.</w:t></w:r><w:r><w:t xml:space="preserve"> (Kommentare)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Jede Stelle des Codes ist nach dem Prinzip clean Code erstellt worden.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="berschrift1"/></w:pPr><w:r><w:t>Verfügbarkeit prüfen</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Jeder kann Benutzer kann die Liste der verfügbaren Räume aufrufen. In dieser Liste kann der Benutzer überprüfen ob ein Raum zu einem bestimmten Zeitpunkt verfügbar ist, oder welche Räume zu einem bestimmten Zeitpunkt verfügbar sind.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Um die Story zu testen, müssen noch einige Räume in die Datenbank eingefügt werden und danach muss das Ganze funktionieren.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Was muss noch gemacht sein: </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Testfälle die im Vorherein geschrieben wurden, sind getestet und ohne Fehler durchgeführt worden.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Dokumentation über die Methoden und Wege wurde geführt und in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> gepuscht. (Kommentare)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Jede Stelle des Codes ist nach dem Prinzip clean Code erstellt worden.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:spacing w:val="-10"/><w:kern w:val="28"/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="Titel"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Definition </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>of</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Done</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Sprint 1</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="berschrift1"/></w:pPr><w:r><w:t>Räume ansehen</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Jeder kann sich eine Liste aller Räume aufrufen und ansehen.</w:t></w:r><w:r><w:br/><w:t>Der User soll in der Liste auf einen Raum klicken können und dann soll eine Beschreibung zu diesem Raum angezeigt werden.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Um die Story zu testen, müssen noch einige Räume in die Datenbank eingefügt werden und danach muss das Ganze funktionieren.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Was muss noch gemacht sein:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Junit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Tests sind geschrieben und alle ohne Fehler ausgeführt worden</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Testfälle die im Vorherein geschrieben wurden, sind getestet und ohne Fehler durchgeführt worden.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Dokumentation über die Methoden und Wege wurde geführt und in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> gepuscht.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Jede Stelle des Codes ist nach dem Prinzip clean Code erstellt worden.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="berschrift1"/></w:pPr><w:r><w:t>Reservationsliste</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Der Hauswart kann sich eine Liste aller Reservationen der nächsten sieben Tage anzeigen lassen.</w:t></w:r><w:r><w:br/><w:t>In den Details werden die Startzeit und die Endzeit sowie die Kontaktdaten des Veranstalters angezeigt, sodass die Räume weiter verplant werden können.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Um die Story zu testen müssen noch einige Reservationen in die Datenbank eingetragen werden, sodass wir das Ganze testen können.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Was muss noch gemacht sein:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Junit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Tests sind geschrieben und alle ohne Fehler ausgeführt worden</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Testfälle die im Vorherein geschrieben wurden, sind getestet und ohne Fehler durchgeführt worden.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Dokumentation über die Methoden und Wege wurde geführt und in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> gepuscht.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Jede Stelle des Codes ist nach dem Prinzip clean Code erstellt worden.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
Write-Host ("Paragraphs after: " + $d.Paragraphs.Count)
